# 036 : Added spacefb driver
#
# The "spacefb" family of games (5 rows) is moved from the "ALL" worksheet
# (rows 1476-1480, global ids 1612-1616) to the end of the
# "Playable (untested)" worksheet (new rows 409-413, renumbered 409-413),
# because they are no longer "not working" but are now considered playable
# (untested). The remaining "ALL" rows shift up by 5, the defined names /
# AutoFilter range shrink accordingly, and the active sheet / selections are
# updated to reflect where the editor last left off.

$wb = $excel.ActiveWorkbook
$wsPlayable = $wb.Worksheets.Item("Playable (untested)")
$wsAll = $wb.Worksheets.Item("ALL")

# ---------------------------------------------------------------------
# 1) Capture the 5 "spacefb" rows (ALL!A1476:M1480) before they move.
# ---------------------------------------------------------------------
$rowsToMove = @(
    @{ A=409; B="spacefb";   C="spacefb.c"; D="Z80"; E="[I8035 ]"; H="1xDAC"; M="Space Firebird (Nintendo)" },
    @{ A=410; B="spacefbg";  C="spacefb.c"; D="Z80"; E="[I8035 ]"; H="1xDAC"; M="Space Firebird (Gremlin)" },
    @{ A=411; B="spacefbb";  C="spacefb.c"; D="Z80"; E="[I8035 ]"; H="1xDAC"; M="Space Firebird (bootleg)" },
    @{ A=412; B="spacebrd";  C="spacefb.c"; D="Z80"; E="[I8035 ]"; H="1xDAC"; M="Space Bird (bootleg)" },
    @{ A=413; B="spacedem";  C="spacefb.c"; D="Z80"; E="[I8035 ]"; H="1xDAC"; M="Space Demon" }
)

# ---------------------------------------------------------------------
# 2) Remove those 5 rows from "ALL" (they currently sit at rows 1476-1480).
#    Everything below shifts up by 5 automatically.
# ---------------------------------------------------------------------
[void]$wsAll.Range("A1476:A1480").EntireRow.Delete()

# ---------------------------------------------------------------------
# 3) Append the moved rows to the bottom of "Playable (untested)"
#    (it previously ended at row 408), renumbering column A to match
#    the new local row number instead of the old global id.
# ---------------------------------------------------------------------
foreach ($row in $rowsToMove) {
    $r = $row.A
    $wsPlayable.Range("A$r").Value2 = $row.A
    $wsPlayable.Range("B$r").Value2 = $row.B
    $wsPlayable.Range("C$r").Value2 = $row.C
    $wsPlayable.Range("D$r").Value2 = $row.D
    $wsPlayable.Range("E$r").Value2 = $row.E
    $wsPlayable.Range("H$r").Value2 = $row.H
    $wsPlayable.Range("M$r").Value2 = $row.M
}

# ---------------------------------------------------------------------
# 4) Shrink the "ALL" defined ranges (_FilterDatabase / LIST) and the
#    AutoFilter to the new extent (M1628 instead of M1633).
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "ALL!_FilterDatabase") {
        $n.RefersTo = "=ALL!`$A`$1:`$M`$1628"
    }
    if ($n.Name -eq "ALL!LIST") {
        $n.RefersTo = "=ALL!`$B`$1:`$M`$1628"
    }
}

$wsAll.AutoFilterMode = $false
[void]$wsAll.Range("A1:M1628").AutoFilter()

# ---------------------------------------------------------------------
# 5) Update selections on both sheets and make "ALL" the active tab,
#    matching where the editor ended up after the move.
# ---------------------------------------------------------------------
[void]$wsPlayable.Range("A414").Select()
[void]$wsAll.Range("E695").Select()
[void]$wsAll.Activate()
